# "estructura inical version final"
# Row 7 of the planning table gets re-arranged / filled in:
#   - B7 becomes "CSS: ventanaModal"   (new task)
#   - C7 becomes "PHP: reservas"       (moved from the old B7)
#   - E7 becomes "PHP: login"          (new task)
# The active selection in the sheet view also moves to C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Write E7 first so the new shared-string entries land in the same order
# ("PHP: login" before "CSS: ventanaModal") as in the target workbook.
$ws.Range("E7").Value = "PHP: login"
$ws.Range("B7").Value = "CSS: ventanaModal"
$ws.Range("C7").Value = "PHP: reservas"

# Move/select C8, matching the updated cursor position saved in the file.
$ws.Range("C8").Select()
